$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for "RM 232" (row 26) and "SC 92" (row 28 before the first delete).
# Deleting row 26 first shifts SC 92 from row 28 to row 27.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# Fix up the individual cell values (which were independently re-imputed)
# in the final (post-delete) row numbering.
$ws.Range("E6").Value = -5.7
$ws.Range("E8").Value = ""
$ws.Range("E19").Value = -6.5
$ws.Range("E21").Value = ""
$ws.Range("E23").Value = -7

$ws.Range("B26").Value = ""
$ws.Range("B27").Value = -20.4
$ws.Range("E27").Value = ""
$ws.Range("B29").Value = ""
$ws.Range("E29").Value = -6.8
